$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.872.14"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").Value = "3.031.64"
$ws.Range("E3").Value = "  -1.30%  "

$ws.Range("E4").Value = "  +0.42%  "

$ws.Range("D5").Value = "586.34"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D6").Value = "151.10"
$ws.Range("E6").Value = "  -1.17%  "

$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("D8").Value = "0.537"
$ws.Range("E8").Value = "  -1.74%  "

$ws.Range("D9").Value = "3.033.62"
$ws.Range("E9").Value = "  -1.38%  "

$ws.Range("E10").Value = "  -2.31%  "

$ws.Range("D11").Value = "5.77"
$ws.Range("E11").Value = "  -1.69%  "

$ws.Range("D12").Value = "0.447"
$ws.Range("E12").Value = "  -3.35%  "

$ws.Range("D13").Value = "0.0000233"
$ws.Range("E13").Value = "  -3.54%  "

$ws.Range("D14").Value = "36.15"
$ws.Range("E14").Value = "  -3.11%  "

$ws.Range("E15").Value = "  +1.32%  "

$ws.Range("D16").Value = "3.531.56"
$ws.Range("E16").Value = "  -1.36%  "

$ws.Range("D17").Value = "7.11"
$ws.Range("E17").Value = "  -1.73%  "

$ws.Range("D18").Value = "62.948.00"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").Value = "3.038.63"
$ws.Range("E19").Value = "  -1.25%  "

$ws.Range("D20").Value = "477.43"
$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("D21").Value = "14.18"
$ws.Range("E21").Value = "  -3.53%  "

$ws.Range("D22").Value = "0.702"
$ws.Range("E22").Value = "  -2.48%  "

$ws.Range("D23").Value = "7.48"
$ws.Range("E23").Value = "  -1.01%  "

$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("D25").Value = "81.80"
$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("D26").Value = "12.68"
$ws.Range("E26").Value = "  -3.76%  "

$ws.Range("D27").Value = "10.60"
$ws.Range("E27").Value = "  +6.31%  "

$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("D29").Value = "7.35"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.70%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.66"
$ws.Range("E31").Value = "  -0.81%  "

$ws.Range("D32").Value = "2.19"
$ws.Range("E32").Value = "  -0.74%  "

$ws.Range("D33").Value = "27.42"
$ws.Range("E33").Value = "  +0.37%  "

$ws.Range("E34").Value = "  -3.40%  "

$ws.Range("D35").Value = "1.06"
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").Value = "0.0₃0811"
$ws.Range("E36").Value = "  -5.05%  "

$ws.Range("D37").Value = "3.25"
$ws.Range("E37").Value = "  -4.15%  "

$ws.Range("D38").Value = "5.90"
$ws.Range("E38").Value = "  -3.98%  "

$ws.Range("D39").Value = "2.21"
$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("D40").Value = "9.22"
$ws.Range("E40").Value = "  -1.23%  "

$ws.Range("D41").Value = "50.30"
$ws.Range("E41").Value = "  -0.17%  "

$ws.Range("D42").Value = "431.74"
$ws.Range("E42").Value = "  -3.47%  "

$ws.Range("D43").Value = "0.286"
$ws.Range("E43").Value = "  -0.31%  "

$ws.Range("D44").Value = "0.113"
$ws.Range("E44").Value = "  +2.05%  "

$ws.Range("D45").Value = "0.0360"
$ws.Range("E45").Value = "  -1.25%  "

$ws.Range("D46").Value = "2.820.60"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("D47").Value = "38.28"
$ws.Range("E47").Value = "  -5.67%  "

$ws.Range("D48").Value = "129.31"
$ws.Range("E48").Value = "  -1.32%  "

$ws.Range("D50").Value = "24.99"
$ws.Range("E50").Value = "  -1.31%  "

$ws.Range("D51").Value = "2.21"
$ws.Range("E51").Value = "  -2.86%  "

